$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cValues = @{
    2  = 7.288845540350142
    3  = -20.40984652067478
    4  = 5.331710924091837
    5  = 11.04982736891558
    6  = 4.748210439985256
    7  = -2.313034291448757
    8  = 3.711391384148
    9  = 3.860244074450181
    10 = 2.370939381494686
    11 = 4.421855465610292
    12 = 3.320585727896552
    13 = 1.782333336406405
    14 = -4.477718018907028
    15 = 6.317691071509768
    16 = 3.79744344971964
    17 = 0.4998689793225486
    18 = 0.1967053802870877
    19 = -2.056549539789942
}

$eValues = @{
    2  = 10.06916370210014
    3  = -33.45158075171641
    4  = 10.84949786623359
    5  = 9.589921161142879
    6  = 7.086193663490992
    7  = -7.020874871669158
    8  = 0.6610955960690834
    9  = 3.254220449867029
    10 = 3.967543131851214
    11 = 5.260364862099642
    12 = -4.098213472638578
    13 = 4.060401000000002
    14 = -11.83522404790002
    15 = 4.613033063261129
    16 = 3.703837953294542
    17 = -0.08221002454066317
    18 = 6.289039804796182
    19 = -0.3884660724497446
}

foreach ($row in 2..19) {
    $ws.Cells.Item($row, 3).Value = $cValues[$row]
    $ws.Cells.Item($row, 5).Value = $eValues[$row]
}
